# Update the "Rating vs Sentiment Matrix" header blocks.
#
# Each of the 8 category blocks (Total, Video_Games, Sports_and_Outdoors,
# Office_Products, Movies_and_TV, Electronics, Clothing_Shoes_and_Jewelry,
# Beauty) has two header rows shaped like:
#
#   row   : B=<category>  C="(down-arrow) Sentiment"                K="Sentiment"
#   row+1 : B="(right-arrow) Score"    C..G = 1..5   J="Score"      K..O = 1..5
#
# The edit swaps the arrow labels (down-arrow now reads "Score", right-arrow
# now reads "Sentiment") and clears out the now-redundant "Sentiment"/"Score"
# labels that used to sit in columns K and J.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$downArrow = [char]::ConvertFromUtf32(0x1F817)
$rightArrow = [char]::ConvertFromUtf32(0x1F816)

$newDownLabel = "$downArrow Score"
$newRightLabel = "$rightArrow Sentiment"

$headerRows = @(30, 38, 46, 54, 62, 70, 78, 86)

foreach ($row in $headerRows) {
    $labelRow = $row
    $dataLabelRow = $row + 1

    # C<row>: "(down-arrow) Sentiment" -> "(down-arrow) Score"
    $ws.Cells.Item($labelRow, 3).Value = $newDownLabel

    # K<row>: "Sentiment" -> cleared
    $ws.Cells.Item($labelRow, 11).Value = ""

    # B<row+1>: "(right-arrow) Score" -> "(right-arrow) Sentiment"
    $ws.Cells.Item($dataLabelRow, 2).Value = $newRightLabel

    # J<row+1>: "Score" -> cleared
    $ws.Cells.Item($dataLabelRow, 10).Value = ""
}
